$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric data columns (D,J,K,L,M,P) for rows 90-219 ---
# This reflects: 2 new price records inserted at rows 90-91 (shifting the rest
# down by two logical positions), with the final two old records now at rows 218-219.

$Dvals = New-Object 'object[,]' 130,1
$Dvals[0,0] = 44895
$Dvals[1,0] = 44895
$Dvals[2,0] = 44781
$Dvals[3,0] = 44671
$Dvals[4,0] = 44893
$Dvals[5,0] = 44383
$Dvals[6,0] = 44343
$Dvals[7,0] = 44221
$Dvals[8,0] = 44350
$Dvals[9,0] = 44341
$Dvals[10,0] = 44376
$Dvals[11,0] = 44312
$Dvals[12,0] = 44453
$Dvals[13,0] = 44389
$Dvals[14,0] = 44217
$Dvals[15,0] = 44663
$Dvals[16,0] = 44420
$Dvals[17,0] = 44620
$Dvals[18,0] = 44316
$Dvals[19,0] = 44567
$Dvals[20,0] = 44386
$Dvals[21,0] = 44308
$Dvals[22,0] = 44358
$Dvals[23,0] = 44648
$Dvals[24,0] = 44382
$Dvals[25,0] = 44330
$Dvals[26,0] = 44670
$Dvals[27,0] = 44426
$Dvals[28,0] = 44629
$Dvals[29,0] = 44448
$Dvals[30,0] = 44532
$Dvals[31,0] = 44399
$Dvals[32,0] = 44461
$Dvals[33,0] = 44719
$Dvals[34,0] = 44285
$Dvals[35,0] = 44284
$Dvals[36,0] = 44634
$Dvals[37,0] = 44680
$Dvals[38,0] = 44209
$Dvals[39,0] = 44812
$Dvals[40,0] = 44582
$Dvals[41,0] = 44315
$Dvals[42,0] = 44357
$Dvals[43,0] = 44859
$Dvals[44,0] = 44477
$Dvals[45,0] = 44615
$Dvals[46,0] = 44342
$Dvals[47,0] = 44609
$Dvals[48,0] = 44300
$Dvals[49,0] = 44356
$Dvals[50,0] = 44879
$Dvals[51,0] = 44665
$Dvals[52,0] = 44525
$Dvals[53,0] = 44340
$Dvals[54,0] = 44776
$Dvals[55,0] = 44195
$Dvals[56,0] = 44838
$Dvals[57,0] = 44658
$Dvals[58,0] = 44813
$Dvals[59,0] = 44313
$Dvals[60,0] = 44428
$Dvals[61,0] = 44809
$Dvals[62,0] = 44333
$Dvals[63,0] = 44412
$Dvals[64,0] = 44784
$Dvals[65,0] = 44286
$Dvals[66,0] = 44245
$Dvals[67,0] = 44411
$Dvals[68,0] = 44397
$Dvals[69,0] = 44188
$Dvals[70,0] = 44274
$Dvals[71,0] = 44720
$Dvals[72,0] = 44165
$Dvals[73,0] = 44627
$Dvals[74,0] = 44868
$Dvals[75,0] = 44336
$Dvals[76,0] = 44476
$Dvals[77,0] = 44260
$Dvals[78,0] = 44301
$Dvals[79,0] = 44398
$Dvals[80,0] = 44767
$Dvals[81,0] = 44214
$Dvals[82,0] = 44638
$Dvals[83,0] = 44371
$Dvals[84,0] = 44385
$Dvals[85,0] = 44736
$Dvals[86,0] = 44679
$Dvals[87,0] = 44881
$Dvals[88,0] = 44847
$Dvals[89,0] = 44372
$Dvals[90,0] = 44889
$Dvals[91,0] = 44811
$Dvals[92,0] = 44782
$Dvals[93,0] = 44441
$Dvals[94,0] = 44649
$Dvals[95,0] = 44196
$Dvals[96,0] = 44425
$Dvals[97,0] = 44771
$Dvals[98,0] = 44272
$Dvals[99,0] = 44581
$Dvals[100,0] = 44406
$Dvals[101,0] = 44161
$Dvals[102,0] = 44446
$Dvals[103,0] = 44704
$Dvals[104,0] = 44685
$Dvals[105,0] = 44370
$Dvals[106,0] = 44211
$Dvals[107,0] = 44294
$Dvals[108,0] = 44413
$Dvals[109,0] = 44522
$Dvals[110,0] = 44662
$Dvals[111,0] = 44203
$Dvals[112,0] = 44568
$Dvals[113,0] = 44495
$Dvals[114,0] = 44650
$Dvals[115,0] = 44298
$Dvals[116,0] = 44763
$Dvals[117,0] = 44455
$Dvals[118,0] = 44753
$Dvals[119,0] = 44299
$Dvals[120,0] = 44622
$Dvals[121,0] = 44610
$Dvals[122,0] = 44162
$Dvals[123,0] = 44636
$Dvals[124,0] = 44414
$Dvals[125,0] = 44837
$Dvals[126,0] = 44364
$Dvals[127,0] = 44400
$Dvals[128,0] = 44218
$Dvals[129,0] = 44417
$ws.Range("D90:D219").Value = $Dvals

$Jvals = New-Object 'object[,]' 130,1
$Jvals[0,0] = 40
$Jvals[1,0] = 20
$Jvals[2,0] = 170
$Jvals[3,0] = 15
$Jvals[4,0] = 40
$Jvals[5,0] = 40
$Jvals[6,0] = 80
$Jvals[7,0] = 55
$Jvals[8,0] = 25
$Jvals[9,0] = 50
$Jvals[10,0] = 15
$Jvals[11,0] = 100
$Jvals[12,0] = 40
$Jvals[13,0] = 35
$Jvals[14,0] = 80
$Jvals[15,0] = 40
$Jvals[16,0] = 15
$Jvals[17,0] = 20
$Jvals[18,0] = 40
$Jvals[19,0] = 50
$Jvals[20,0] = 30
$Jvals[21,0] = 35
$Jvals[22,0] = 50
$Jvals[23,0] = 130
$Jvals[24,0] = 30
$Jvals[25,0] = 50
$Jvals[26,0] = 15
$Jvals[27,0] = 50
$Jvals[28,0] = 15
$Jvals[29,0] = 55
$Jvals[30,0] = 30
$Jvals[31,0] = 100
$Jvals[32,0] = 30
$Jvals[33,0] = 50
$Jvals[34,0] = 80
$Jvals[35,0] = 100
$Jvals[36,0] = 70
$Jvals[37,0] = 80
$Jvals[38,0] = 100
$Jvals[39,0] = 160
$Jvals[40,0] = 40
$Jvals[41,0] = 30
$Jvals[42,0] = 50
$Jvals[43,0] = 15
$Jvals[44,0] = 30
$Jvals[45,0] = 10
$Jvals[46,0] = 40
$Jvals[47,0] = 20
$Jvals[48,0] = 30
$Jvals[49,0] = 60
$Jvals[50,0] = 90
$Jvals[51,0] = 60
$Jvals[52,0] = 20
$Jvals[53,0] = 80
$Jvals[54,0] = 30
$Jvals[55,0] = 75
$Jvals[56,0] = 20
$Jvals[57,0] = 65
$Jvals[58,0] = 40
$Jvals[59,0] = 50
$Jvals[60,0] = 120
$Jvals[61,0] = 120
$Jvals[62,0] = 15
$Jvals[63,0] = 30
$Jvals[64,0] = 100
$Jvals[65,0] = 80
$Jvals[66,0] = 45
$Jvals[67,0] = 20
$Jvals[68,0] = 30
$Jvals[69,0] = 40
$Jvals[70,0] = 50
$Jvals[71,0] = 50
$Jvals[72,0] = 65
$Jvals[73,0] = 35
$Jvals[74,0] = 200
$Jvals[75,0] = 50
$Jvals[76,0] = 50
$Jvals[77,0] = 40
$Jvals[78,0] = 80
$Jvals[79,0] = 80
$Jvals[80,0] = 50
$Jvals[81,0] = 50
$Jvals[82,0] = 50
$Jvals[83,0] = 80
$Jvals[84,0] = 80
$Jvals[85,0] = 40
$Jvals[86,0] = 50
$Jvals[87,0] = 30
$Jvals[88,0] = 25
$Jvals[89,0] = 20
$Jvals[90,0] = 6
$Jvals[91,0] = 100
$Jvals[92,0] = 140
$Jvals[93,0] = 30
$Jvals[94,0] = 50
$Jvals[95,0] = 20
$Jvals[96,0] = 30
$Jvals[97,0] = 40
$Jvals[98,0] = 40
$Jvals[99,0] = 40
$Jvals[100,0] = 20
$Jvals[101,0] = 130
$Jvals[102,0] = 25
$Jvals[103,0] = 40
$Jvals[104,0] = 50
$Jvals[105,0] = 20
$Jvals[106,0] = 35
$Jvals[107,0] = 15
$Jvals[108,0] = 100
$Jvals[109,0] = 60
$Jvals[110,0] = 100
$Jvals[111,0] = 50
$Jvals[112,0] = 40
$Jvals[113,0] = 30
$Jvals[114,0] = 60
$Jvals[115,0] = 80
$Jvals[116,0] = 25
$Jvals[117,0] = 50
$Jvals[118,0] = 80
$Jvals[119,0] = 40
$Jvals[120,0] = 40
$Jvals[121,0] = 40
$Jvals[122,0] = 130
$Jvals[123,0] = 40
$Jvals[124,0] = 80
$Jvals[125,0] = 80
$Jvals[126,0] = 40
$Jvals[127,0] = 40
$Jvals[128,0] = 60
$Jvals[129,0] = 45
$ws.Range("J90:J219").Value = $Jvals

$Kvals = New-Object 'object[,]' 130,1
$Kvals[0,0] = 20000
$Kvals[1,0] = 20000
$Kvals[2,0] = 16000
$Kvals[3,0] = 20000
$Kvals[4,0] = 20000
$Kvals[5,0] = 25000
$Kvals[6,0] = 30000
$Kvals[7,0] = 35000
$Kvals[8,0] = 28000
$Kvals[9,0] = 30000
$Kvals[10,0] = 25000
$Kvals[11,0] = 35000
$Kvals[12,0] = 20000
$Kvals[13,0] = 20000
$Kvals[14,0] = 35000
$Kvals[15,0] = 25000
$Kvals[16,0] = 25000
$Kvals[17,0] = 26000
$Kvals[18,0] = 35000
$Kvals[19,0] = 20000
$Kvals[20,0] = 25000
$Kvals[21,0] = 35000
$Kvals[22,0] = 35000
$Kvals[23,0] = 20000
$Kvals[24,0] = 25000
$Kvals[25,0] = 28000
$Kvals[26,0] = 20000
$Kvals[27,0] = 30000
$Kvals[28,0] = 25000
$Kvals[29,0] = 20000
$Kvals[30,0] = 20000
$Kvals[31,0] = 20000
$Kvals[32,0] = 20000
$Kvals[33,0] = 20000
$Kvals[34,0] = 32000
$Kvals[35,0] = 32000
$Kvals[36,0] = 22000
$Kvals[37,0] = 20000
$Kvals[38,0] = 30000
$Kvals[39,0] = 20000
$Kvals[40,0] = 20000
$Kvals[41,0] = 35000
$Kvals[42,0] = 35000
$Kvals[43,0] = 20000
$Kvals[44,0] = 20000
$Kvals[45,0] = 26000
$Kvals[46,0] = 30000
$Kvals[47,0] = 26000
$Kvals[48,0] = 35000
$Kvals[49,0] = 32000
$Kvals[50,0] = 18000
$Kvals[51,0] = 15000
$Kvals[52,0] = 20000
$Kvals[53,0] = 25000
$Kvals[54,0] = 20000
$Kvals[55,0] = 30000
$Kvals[56,0] = 20000
$Kvals[57,0] = 25000
$Kvals[58,0] = 20000
$Kvals[59,0] = 35000
$Kvals[60,0] = 20000
$Kvals[61,0] = 20000
$Kvals[62,0] = 35000
$Kvals[63,0] = 20000
$Kvals[64,0] = 16000
$Kvals[65,0] = 32000
$Kvals[66,0] = 35000
$Kvals[67,0] = 20000
$Kvals[68,0] = 21000
$Kvals[69,0] = 33000
$Kvals[70,0] = 40000
$Kvals[71,0] = 20000
$Kvals[72,0] = 35000
$Kvals[73,0] = 25000
$Kvals[74,0] = 18000
$Kvals[75,0] = 28000
$Kvals[76,0] = 20000
$Kvals[77,0] = 40000
$Kvals[78,0] = 35000
$Kvals[79,0] = 21000
$Kvals[80,0] = 20000
$Kvals[81,0] = 35000
$Kvals[82,0] = 20000
$Kvals[83,0] = 25000
$Kvals[84,0] = 25000
$Kvals[85,0] = 20000
$Kvals[86,0] = 20000
$Kvals[87,0] = 20000
$Kvals[88,0] = 20000
$Kvals[89,0] = 25000
$Kvals[90,0] = 20000
$Kvals[91,0] = 18000
$Kvals[92,0] = 16000
$Kvals[93,0] = 25000
$Kvals[94,0] = 20000
$Kvals[95,0] = 30000
$Kvals[96,0] = 25000
$Kvals[97,0] = 20000
$Kvals[98,0] = 40000
$Kvals[99,0] = 20000
$Kvals[100,0] = 25000
$Kvals[101,0] = 33000
$Kvals[102,0] = 20000
$Kvals[103,0] = 20000
$Kvals[104,0] = 20000
$Kvals[105,0] = 25000
$Kvals[106,0] = 35000
$Kvals[107,0] = 35000
$Kvals[108,0] = 25000
$Kvals[109,0] = 20000
$Kvals[110,0] = 25000
$Kvals[111,0] = 30000
$Kvals[112,0] = 20000
$Kvals[113,0] = 25000
$Kvals[114,0] = 20000
$Kvals[115,0] = 35000
$Kvals[116,0] = 20000
$Kvals[117,0] = 20000
$Kvals[118,0] = 20000
$Kvals[119,0] = 35000
$Kvals[120,0] = 26000
$Kvals[121,0] = 25000
$Kvals[122,0] = 33000
$Kvals[123,0] = 25000
$Kvals[124,0] = 20000
$Kvals[125,0] = 20000
$Kvals[126,0] = 25000
$Kvals[127,0] = 20000
$Kvals[128,0] = 35000
$Kvals[129,0] = 25000
$ws.Range("K90:K219").Value = $Kvals

$Lvals = New-Object 'object[,]' 130,1
$Lvals[0,0] = 20000
$Lvals[1,0] = 20000
$Lvals[2,0] = 18000
$Lvals[3,0] = 20000
$Lvals[4,0] = 20000
$Lvals[5,0] = 25000
$Lvals[6,0] = 30000
$Lvals[7,0] = 35000
$Lvals[8,0] = 28000
$Lvals[9,0] = 32000
$Lvals[10,0] = 25000
$Lvals[11,0] = 35000
$Lvals[12,0] = 25000
$Lvals[13,0] = 20000
$Lvals[14,0] = 35000
$Lvals[15,0] = 25000
$Lvals[16,0] = 25000
$Lvals[17,0] = 26000
$Lvals[18,0] = 35000
$Lvals[19,0] = 20000
$Lvals[20,0] = 25000
$Lvals[21,0] = 35000
$Lvals[22,0] = 35000
$Lvals[23,0] = 25000
$Lvals[24,0] = 25000
$Lvals[25,0] = 28000
$Lvals[26,0] = 20000
$Lvals[27,0] = 30000
$Lvals[28,0] = 25000
$Lvals[29,0] = 20000
$Lvals[30,0] = 20000
$Lvals[31,0] = 20000
$Lvals[32,0] = 20000
$Lvals[33,0] = 20000
$Lvals[34,0] = 32000
$Lvals[35,0] = 32000
$Lvals[36,0] = 25000
$Lvals[37,0] = 20000
$Lvals[38,0] = 30000
$Lvals[39,0] = 20000
$Lvals[40,0] = 20000
$Lvals[41,0] = 35000
$Lvals[42,0] = 35000
$Lvals[43,0] = 20000
$Lvals[44,0] = 20000
$Lvals[45,0] = 26000
$Lvals[46,0] = 30000
$Lvals[47,0] = 26000
$Lvals[48,0] = 35000
$Lvals[49,0] = 35000
$Lvals[50,0] = 20000
$Lvals[51,0] = 25000
$Lvals[52,0] = 20000
$Lvals[53,0] = 25000
$Lvals[54,0] = 20000
$Lvals[55,0] = 30000
$Lvals[56,0] = 20000
$Lvals[57,0] = 25000
$Lvals[58,0] = 20000
$Lvals[59,0] = 35000
$Lvals[60,0] = 25000
$Lvals[61,0] = 20000
$Lvals[62,0] = 35000
$Lvals[63,0] = 20000
$Lvals[64,0] = 16000
$Lvals[65,0] = 32000
$Lvals[66,0] = 35000
$Lvals[67,0] = 20000
$Lvals[68,0] = 21000
$Lvals[69,0] = 33000
$Lvals[70,0] = 40000
$Lvals[71,0] = 20000
$Lvals[72,0] = 35000
$Lvals[73,0] = 26000
$Lvals[74,0] = 20000
$Lvals[75,0] = 28000
$Lvals[76,0] = 20000
$Lvals[77,0] = 40000
$Lvals[78,0] = 35000
$Lvals[79,0] = 21000
$Lvals[80,0] = 20000
$Lvals[81,0] = 35000
$Lvals[82,0] = 20000
$Lvals[83,0] = 25000
$Lvals[84,0] = 25000
$Lvals[85,0] = 20000
$Lvals[86,0] = 20000
$Lvals[87,0] = 20000
$Lvals[88,0] = 20000
$Lvals[89,0] = 25000
$Lvals[90,0] = 20000
$Lvals[91,0] = 20000
$Lvals[92,0] = 20000
$Lvals[93,0] = 25000
$Lvals[94,0] = 20000
$Lvals[95,0] = 30000
$Lvals[96,0] = 25000
$Lvals[97,0] = 20000
$Lvals[98,0] = 40000
$Lvals[99,0] = 20000
$Lvals[100,0] = 25000
$Lvals[101,0] = 35000
$Lvals[102,0] = 20000
$Lvals[103,0] = 20000
$Lvals[104,0] = 20000
$Lvals[105,0] = 25000
$Lvals[106,0] = 35000
$Lvals[107,0] = 35000
$Lvals[108,0] = 25000
$Lvals[109,0] = 20000
$Lvals[110,0] = 25000
$Lvals[111,0] = 30000
$Lvals[112,0] = 20000
$Lvals[113,0] = 25000
$Lvals[114,0] = 25000
$Lvals[115,0] = 35000
$Lvals[116,0] = 20000
$Lvals[117,0] = 20000
$Lvals[118,0] = 20000
$Lvals[119,0] = 35000
$Lvals[120,0] = 26000
$Lvals[121,0] = 26000
$Lvals[122,0] = 35000
$Lvals[123,0] = 25000
$Lvals[124,0] = 20000
$Lvals[125,0] = 20000
$Lvals[126,0] = 25000
$Lvals[127,0] = 20000
$Lvals[128,0] = 35000
$Lvals[129,0] = 26000
$ws.Range("L90:L219").Value = $Lvals

$Mvals = New-Object 'object[,]' 130,1
$Mvals[0,0] = 20000
$Mvals[1,0] = 20000
$Mvals[2,0] = 16941
$Mvals[3,0] = 20000
$Mvals[4,0] = 20000
$Mvals[5,0] = 25000
$Mvals[6,0] = 30000
$Mvals[7,0] = 35000
$Mvals[8,0] = 28000
$Mvals[9,0] = 30800
$Mvals[10,0] = 25000
$Mvals[11,0] = 35000
$Mvals[12,0] = 22500
$Mvals[13,0] = 20000
$Mvals[14,0] = 35000
$Mvals[15,0] = 25000
$Mvals[16,0] = 25000
$Mvals[17,0] = 26000
$Mvals[18,0] = 35000
$Mvals[19,0] = 20000
$Mvals[20,0] = 25000
$Mvals[21,0] = 35000
$Mvals[22,0] = 35000
$Mvals[23,0] = 21923
$Mvals[24,0] = 25000
$Mvals[25,0] = 28000
$Mvals[26,0] = 20000
$Mvals[27,0] = 30000
$Mvals[28,0] = 25000
$Mvals[29,0] = 20000
$Mvals[30,0] = 20000
$Mvals[31,0] = 20000
$Mvals[32,0] = 20000
$Mvals[33,0] = 20000
$Mvals[34,0] = 32000
$Mvals[35,0] = 32000
$Mvals[36,0] = 23714
$Mvals[37,0] = 20000
$Mvals[38,0] = 30000
$Mvals[39,0] = 20000
$Mvals[40,0] = 20000
$Mvals[41,0] = 35000
$Mvals[42,0] = 35000
$Mvals[43,0] = 20000
$Mvals[44,0] = 20000
$Mvals[45,0] = 26000
$Mvals[46,0] = 30000
$Mvals[47,0] = 26000
$Mvals[48,0] = 35000
$Mvals[49,0] = 34000
$Mvals[50,0] = 18889
$Mvals[51,0] = 21667
$Mvals[52,0] = 20000
$Mvals[53,0] = 25000
$Mvals[54,0] = 20000
$Mvals[55,0] = 30000
$Mvals[56,0] = 20000
$Mvals[57,0] = 25000
$Mvals[58,0] = 20000
$Mvals[59,0] = 35000
$Mvals[60,0] = 22083
$Mvals[61,0] = 20000
$Mvals[62,0] = 35000
$Mvals[63,0] = 20000
$Mvals[64,0] = 16000
$Mvals[65,0] = 32000
$Mvals[66,0] = 35000
$Mvals[67,0] = 20000
$Mvals[68,0] = 21000
$Mvals[69,0] = 33000
$Mvals[70,0] = 40000
$Mvals[71,0] = 20000
$Mvals[72,0] = 35000
$Mvals[73,0] = 25571
$Mvals[74,0] = 19000
$Mvals[75,0] = 28000
$Mvals[76,0] = 20000
$Mvals[77,0] = 40000
$Mvals[78,0] = 35000
$Mvals[79,0] = 21000
$Mvals[80,0] = 20000
$Mvals[81,0] = 35000
$Mvals[82,0] = 20000
$Mvals[83,0] = 25000
$Mvals[84,0] = 25000
$Mvals[85,0] = 20000
$Mvals[86,0] = 20000
$Mvals[87,0] = 20000
$Mvals[88,0] = 20000
$Mvals[89,0] = 25000
$Mvals[90,0] = 20000
$Mvals[91,0] = 19000
$Mvals[92,0] = 17429
$Mvals[93,0] = 25000
$Mvals[94,0] = 20000
$Mvals[95,0] = 30000
$Mvals[96,0] = 25000
$Mvals[97,0] = 20000
$Mvals[98,0] = 40000
$Mvals[99,0] = 20000
$Mvals[100,0] = 25000
$Mvals[101,0] = 33769
$Mvals[102,0] = 20000
$Mvals[103,0] = 20000
$Mvals[104,0] = 20000
$Mvals[105,0] = 25000
$Mvals[106,0] = 35000
$Mvals[107,0] = 35000
$Mvals[108,0] = 25000
$Mvals[109,0] = 20000
$Mvals[110,0] = 25000
$Mvals[111,0] = 30000
$Mvals[112,0] = 20000
$Mvals[113,0] = 25000
$Mvals[114,0] = 21667
$Mvals[115,0] = 35000
$Mvals[116,0] = 20000
$Mvals[117,0] = 20000
$Mvals[118,0] = 20000
$Mvals[119,0] = 35000
$Mvals[120,0] = 26000
$Mvals[121,0] = 25500
$Mvals[122,0] = 33769
$Mvals[123,0] = 25000
$Mvals[124,0] = 20000
$Mvals[125,0] = 20000
$Mvals[126,0] = 25000
$Mvals[127,0] = 20000
$Mvals[128,0] = 35000
$Mvals[129,0] = 25556
$ws.Range("M90:M219").Value = $Mvals

$Pvals = New-Object 'object[,]' 130,1
$Pvals[0,0] = 1538
$Pvals[1,0] = 1538
$Pvals[2,0] = 1303
$Pvals[3,0] = 1538
$Pvals[4,0] = 1538
$Pvals[5,0] = 1923
$Pvals[6,0] = 2308
$Pvals[7,0] = 2692
$Pvals[8,0] = 2154
$Pvals[9,0] = 2369
$Pvals[10,0] = 1923
$Pvals[11,0] = 2692
$Pvals[12,0] = 1731
$Pvals[13,0] = 1538
$Pvals[14,0] = 2692
$Pvals[15,0] = 1923
$Pvals[16,0] = 1923
$Pvals[17,0] = 2000
$Pvals[18,0] = 2692
$Pvals[19,0] = 1538
$Pvals[20,0] = 1923
$Pvals[21,0] = 2692
$Pvals[22,0] = 2692
$Pvals[23,0] = 1686
$Pvals[24,0] = 1923
$Pvals[25,0] = 2154
$Pvals[26,0] = 1538
$Pvals[27,0] = 2308
$Pvals[28,0] = 1923
$Pvals[29,0] = 1538
$Pvals[30,0] = 1538
$Pvals[31,0] = 1538
$Pvals[32,0] = 1538
$Pvals[33,0] = 1538
$Pvals[34,0] = 2462
$Pvals[35,0] = 2462
$Pvals[36,0] = 1824
$Pvals[37,0] = 1538
$Pvals[38,0] = 2308
$Pvals[39,0] = 1538
$Pvals[40,0] = 1538
$Pvals[41,0] = 2692
$Pvals[42,0] = 2692
$Pvals[43,0] = 1538
$Pvals[44,0] = 1538
$Pvals[45,0] = 2000
$Pvals[46,0] = 2308
$Pvals[47,0] = 2000
$Pvals[48,0] = 2692
$Pvals[49,0] = 2615
$Pvals[50,0] = 1453
$Pvals[51,0] = 1667
$Pvals[52,0] = 1538
$Pvals[53,0] = 1923
$Pvals[54,0] = 1538
$Pvals[55,0] = 2308
$Pvals[56,0] = 1538
$Pvals[57,0] = 1923
$Pvals[58,0] = 1538
$Pvals[59,0] = 2692
$Pvals[60,0] = 1699
$Pvals[61,0] = 1538
$Pvals[62,0] = 2692
$Pvals[63,0] = 1538
$Pvals[64,0] = 1231
$Pvals[65,0] = 2462
$Pvals[66,0] = 2692
$Pvals[67,0] = 1538
$Pvals[68,0] = 1615
$Pvals[69,0] = 2538
$Pvals[70,0] = 3077
$Pvals[71,0] = 1538
$Pvals[72,0] = 2692
$Pvals[73,0] = 1967
$Pvals[74,0] = 1462
$Pvals[75,0] = 2154
$Pvals[76,0] = 1538
$Pvals[77,0] = 3077
$Pvals[78,0] = 2692
$Pvals[79,0] = 1615
$Pvals[80,0] = 1538
$Pvals[81,0] = 2692
$Pvals[82,0] = 1538
$Pvals[83,0] = 1923
$Pvals[84,0] = 1923
$Pvals[85,0] = 1538
$Pvals[86,0] = 1538
$Pvals[87,0] = 1538
$Pvals[88,0] = 1538
$Pvals[89,0] = 1923
$Pvals[90,0] = 1538
$Pvals[91,0] = 1462
$Pvals[92,0] = 1341
$Pvals[93,0] = 1923
$Pvals[94,0] = 1538
$Pvals[95,0] = 2308
$Pvals[96,0] = 1923
$Pvals[97,0] = 1538
$Pvals[98,0] = 3077
$Pvals[99,0] = 1538
$Pvals[100,0] = 1923
$Pvals[101,0] = 2598
$Pvals[102,0] = 1538
$Pvals[103,0] = 1538
$Pvals[104,0] = 1538
$Pvals[105,0] = 1923
$Pvals[106,0] = 2692
$Pvals[107,0] = 2692
$Pvals[108,0] = 1923
$Pvals[109,0] = 1538
$Pvals[110,0] = 1923
$Pvals[111,0] = 2308
$Pvals[112,0] = 1538
$Pvals[113,0] = 1923
$Pvals[114,0] = 1667
$Pvals[115,0] = 2692
$Pvals[116,0] = 1538
$Pvals[117,0] = 1538
$Pvals[118,0] = 1538
$Pvals[119,0] = 2692
$Pvals[120,0] = 2000
$Pvals[121,0] = 1962
$Pvals[122,0] = 2598
$Pvals[123,0] = 1923
$Pvals[124,0] = 1538
$Pvals[125,0] = 1538
$Pvals[126,0] = 1923
$Pvals[127,0] = 1538
$Pvals[128,0] = 2692
$Pvals[129,0] = 1966
$ws.Range("P90:P219").Value = $Pvals

# --- Populate the constant (unchanging) columns for newly appended rows 218:219 ---
$ws.Range("A218:A219").Value = 10
$ws.Range("B218:B219").Value = "Vega Modelo de Temuco"
$ws.Range("C218:C219").Value = "La Araucanía"
$ws.Range("E218:E219").Value = 9
$ws.Range("F218:F219").Value = 100114007
$ws.Range("G218:G219").Value = "Jengibre"
$ws.Range("H218:H219").Value = "Sin especificar"
$ws.Range("I218:I219").Value = "Primera"
$ws.Range("N218:N219").Value = "$/caja 13 kilos"
$ws.Range("O218:O219").Value = "Perú"
$ws.Range("Q218:Q219").Value = 13
$ws.Range("R218:R219").Value = "Hortaliza"

# --- Apply the custom date/time number format to column D for the full updated range ---
$ws.Range("D90:D219").NumberFormat = "YYYY-MM-DD HH:MM:SS"